# Updates cryptocurrency Price (D) and Volume(1h) (E) columns with
# refreshed values, matching the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Both columns hold plain text (prices use "." as both a thousands and
# decimal separator, e.g. "30.448.11", and volumes are space-padded
# percentage strings, e.g. "  -0.95%  "). Values that otherwise look like
# an ordinary decimal number (e.g. "329.67") are written with a leading
# apostrophe so Excel keeps them as text instead of auto-converting them
# to a Number (which would also silently drop trailing zeros).

$ws.Range("D2").Value = "30.448.11"
$ws.Range("E2").Value = "  -0.95%  "
$ws.Range("D3").Value = "2.091.77"
$ws.Range("E3").Value = "  -1.23%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "'329.67"
$ws.Range("E5").Value = "  -1.30%  "
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("D7").Value = "'0.5199"
$ws.Range("E7").Value = "  -0.60%  "
$ws.Range("D8").Value = "'0.4362"
$ws.Range("E8").Value = "  -1.42%  "
$ws.Range("D9").Value = "'54.03"
$ws.Range("E9").Value = "  +16.58%  "
$ws.Range("D10").Value = "'0.08833"
$ws.Range("E10").Value = "  -2.86%  "
$ws.Range("E11").Value = "  -2.78%  "
$ws.Range("E12").Value = "  -4.18%  "
$ws.Range("D13").Value = "2.075.00"
$ws.Range("E13").Value = "  -1.76%  "
$ws.Range("D14").Value = "'6.686"
$ws.Range("E14").Value = "  -1.86%  "
$ws.Range("D15").Value = "'7.670"
$ws.Range("E15").Value = "  -1.24%  "
$ws.Range("E16").Value = "  -2.89%  "
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").Value = "'0.00001119"
$ws.Range("D19").Value = "'0.06588"
$ws.Range("E19").Value = "  -1.05%  "
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("D22").Value = "'6.267"
$ws.Range("E22").Value = "  -2.26%  "
$ws.Range("D23").Value = "30.479.80"
$ws.Range("E23").Value = "  -1.18%  "
$ws.Range("D24").Value = "'12.23"
$ws.Range("E24").Value = "  +1.24%  "
$ws.Range("D26").Value = "2.331.60"
$ws.Range("E26").Value = "  -1.25%  "
$ws.Range("E27").Value = "  -3.42%  "
$ws.Range("D28").Value = "'2.561"
$ws.Range("E28").Value = "  -0.29%  "
$ws.Range("D29").Value = "'162.69"
$ws.Range("E29").Value = "  -0.58%  "
$ws.Range("D30").Value = "'131.51"
$ws.Range("E30").Value = "  -1.82%  "
$ws.Range("D31").Value = "'1.184"
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("E33").Value = "  +7.10%  "
$ws.Range("D34").Value = "'6.160"
$ws.Range("E34").Value = "  -1.55%  "
$ws.Range("D35").Value = "'3.905"
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("D36").Value = "'10.07"
$ws.Range("E36").Value = "  +4.55%  "
$ws.Range("E37").Value = "  -0.66%  "
$ws.Range("D38").Value = "'0.06794"
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("D39").Value = "'5.446"
$ws.Range("E39").Value = "  -2.82%  "
$ws.Range("E40").Value = "  -1.27%  "
$ws.Range("D41").Value = "'0.2250"
$ws.Range("E41").Value = "  -1.00%  "
$ws.Range("D42").Value = "'0.6861"
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("D43").Value = "'1.255"
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("D45").Value = "'0.6330"
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("E46").Value = "  -3.47%  "
$ws.Range("D47").Value = "'2.191"
$ws.Range("E47").Value = "  -3.00%  "
$ws.Range("D48").Value = "'3.626"
$ws.Range("E48").Value = "  -1.36%  "
$ws.Range("E49").Value = "  +7.07%  "
$ws.Range("D50").Value = "'1.242"
$ws.Range("E50").Value = "  -3.89%  "
$ws.Range("E51").Value = "  -1.97%  "
